# Update "想去人数" (want-to-go count) values in the F column on the
# "展览" and "全部类型" worksheets to reflect the latest generated data.

$wb = $excel.ActiveWorkbook

$sheetUpdates = @{
    "展览"   = @{ "F2" = 1234; "F5" = 12526; "F6" = 70; "F9" = 7; "F10" = 12403; "F13" = 4820; "F16" = 425; "F17" = 107; "F19" = 11 }
    "全部类型" = @{ "F2" = 1234; "F7" = 12526; "F8" = 70; "F11" = 7; "F12" = 12403; "F15" = 4820; "F18" = 425; "F19" = 107; "F21" = 11 }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $sheetUpdates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
